$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52:94 down to 53:95.
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new weekly record
# (same data as the former last record, row 93 -> now 94, but with a new date).
$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = 44810
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112021
$ws.Cells.Item(52, 7).Value = "Ají"
$ws.Cells.Item(52, 8).Value = "Inferno"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 15000
$ws.Cells.Item(52, 12).Value = 16000
$ws.Cells.Item(52, 13).Value = 15500
$ws.Cells.Item(52, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 1033
$ws.Cells.Item(52, 17).Value = 15
$ws.Cells.Item(52, 18).Value = "Hortaliza"
